# Apply the cyclic rotation of data rows 2, 3, 4:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# Only the columns that actually differ between the three rows are touched
# (A, B, C, E, F, G, H, I, J, L, P, Q, R, Y, AA, AX); columns that are
# identical across the three source rows (D, K, M, N, S, T, U, V, W, Z, AB,
# AD, AE, AG, AT, AW, AY) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel silently
# reinterpreting it as a number/date (which it would do for strings like
# "1999-04-28" or "2"). Prefixing with an apostrophe forces text entry,
# then resetting the style back to Normal drops the quote-prefix /
# number-format flag that the apostrophe trick leaves behind, so the
# written cell ends up plain text with no extra styling -- matching how
# the source file stores these as bare inline strings.
function Set-TextCell {
    param($Sheet, $Row, $Col, [string]$Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

function Set-NumCell {
    param($Sheet, $Row, $Col, $Number)
    $Sheet.Cells.Item($Row, $Col).Value = $Number
}

# ---- Row 2 (becomes the old row 3 content) ----
Set-NumCell  $ws 2 1  68175904
Set-NumCell  $ws 2 2  57585
Set-TextCell $ws 2 3  "Ovaliderad"
Set-NumCell  $ws 2 5  208242
Set-TextCell $ws 2 6  "Mindre vattensalamander"
Set-TextCell $ws 2 7  "Lissotriton vulgaris"
Set-TextCell $ws 2 8  "(Linnaeus, 1758)"
Set-TextCell $ws 2 9  "2"
Set-TextCell $ws 2 10 "ex."
Set-TextCell $ws 2 12 "hona"
Set-TextCell $ws 2 16 "Sjövik, Sommen, Sm"
Set-NumCell  $ws 2 17 503498.5757228022
Set-NumCell  $ws 2 18 6429815.746484536
Set-TextCell $ws 2 25 "1999-04-28"
Set-TextCell $ws 2 27 "1999-04-28"
Set-TextCell $ws 2 50 "Josefine Gustafsson"

# ---- Row 3 (becomes the old row 4 content) ----
Set-NumCell  $ws 3 1  68175906
Set-NumCell  $ws 3 2  57587
Set-TextCell $ws 3 3  "Ovaliderad"
Set-NumCell  $ws 3 5  100141
Set-TextCell $ws 3 6  "Större vattensalamander"
Set-TextCell $ws 3 7  "Triturus cristatus"
Set-TextCell $ws 3 8  "(Laurenti, 1768)"
Set-TextCell $ws 3 9  "1"
Set-TextCell $ws 3 10 "ex."
Set-TextCell $ws 3 12 "hona"
Set-TextCell $ws 3 16 "Sjövik, Sommen, Sm"
Set-NumCell  $ws 3 17 503498.5757228022
Set-NumCell  $ws 3 18 6429815.746484536
Set-TextCell $ws 3 25 "1999-04-28"
Set-TextCell $ws 3 27 "1999-04-28"
Set-TextCell $ws 3 50 "Josefine Gustafsson"

# ---- Row 4 (becomes the old row 2 content) ----
Set-NumCell  $ws 4 1  68176169
Set-NumCell  $ws 4 2  57587
Set-TextCell $ws 4 3  "Behöver inte valideras"
Set-NumCell  $ws 4 5  100141
Set-TextCell $ws 4 6  "Större vattensalamander"
Set-TextCell $ws 4 7  "Triturus cristatus"
Set-TextCell $ws 4 8  "(Laurenti, 1768)"
Set-TextCell $ws 4 9  ""
Set-TextCell $ws 4 10 ""
Set-TextCell $ws 4 12 ""
Set-TextCell $ws 4 16 "Sjövik, Sm"
Set-NumCell  $ws 4 17 503526.6896539551
Set-NumCell  $ws 4 18 6429839.084042171
Set-TextCell $ws 4 25 "2005-01-01"
Set-TextCell $ws 4 27 "2005-12-31"
Set-TextCell $ws 4 50 "Anna Isaksson"
